$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (shifts old N,O,P -> O,P,Q)
$ws.Columns("N").Insert() | Out-Null
$ws.Columns("N").ColumnWidth = 9.8

# Make "Repayment schedule" the active sheet/tab and set its selection
$ws.Activate() | Out-Null
$ws.Range("G18").Select() | Out-Null

Write-Host "Done"
